$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows (74-76) to the termData table for a new "New Vocab" /
# "newvocab" term list with terms "Term 1", "Term 2", "Term 3".

# Row 74: New Vocab / newvocab / Term 1
$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "create"
$ws.Range("D74").Value = "New Vocab"
$ws.Range("E74").Value = "newvocab"
$ws.Range("F74").Value = "Term 1"
$ws.Range("L74").Value = "Term 1"
$ws.Range("C74").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,E74,L74)"
$ws.Range("K74").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,C74,A74)"

# Row 75: New Vocab / newvocab / Term 2
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "create"
$ws.Range("D75").Value = "New Vocab"
$ws.Range("E75").Value = "newvocab"
$ws.Range("F75").Value = "Term 2"
$ws.Range("L75").Value = "Term 2"
$ws.Range("C75").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,E75,L75)"
$ws.Range("K75").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,C75,A75)"

# Row 76: New Vocab / newvocab / Term 3
$ws.Range("A76").Value = 4
$ws.Range("B76").Value = "create"
$ws.Range("D76").Value = "New Vocab"
$ws.Range("E76").Value = "newvocab"
$ws.Range("F76").Value = "Term 3"
$ws.Range("L76").Value = "Term 3"
$ws.Range("C76").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,E76,L76)"
$ws.Range("K76").Formula = "=_xlfn.TEXTJOIN("" "",TRUE,C76,A76)"

# The "sort-dedupe" column (K) picks up a distinct number-format style for
# these new rows (mirrors the workbook author applying a format to the new
# rows before filling them in).
$ws.Range("K74:K76").NumberFormat = "General"

# Grow the table so the new rows become part of it (and the autofilter
# range grows along with it).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L76"))

# Reflect the author's final selection/scroll position on the sheet.
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L75:L76").Select()

Write-Host "Added New Vocab term rows (74-76) and resized termData table"
